$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 44318
$ws.Range("A6").NumberFormat = "d-mmm"

$ws.Range("B6").Value = 0.083333333333333329
$ws.Range("B6").NumberFormat = "h:mm"

$ws.Range("C6").Value = "fonctionnement du jeu"

$ws.Range("C7").Select()
